$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that needs to move from
# 46075 (2026-02-22) to 46076 (2026-02-23) for every data row (2-37).
for ($row = 2; $row -le 37; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
